# milestones.xlsx - ch7: fill in "資料庫建置" milestone rows (15-19)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D, top to bottom ------------------------------------------------
# D15 already held the placeholder "---"; replace it with the real milestone
# item name.
$ws.Range("D15").Value = "資料庫建置"

$ws.Range("D16").Value = "Rails 建置"
# "建置" portion of D16 gets its own run formatting (sz 12, white, 宋体)
$d16chars = $ws.Range("D16").Characters(7, 2)
$d16chars.Font.Name = "宋体"
$d16chars.Font.Size = 12
$d16chars.Font.Color = 16777215

$ws.Range("D17").Value = "使用者登入登出頁面"
$ws.Range("D18").Value = "使用者登入"
$ws.Range("D19").Value = "個人資訊編輯"

# --- Column E, rows 17-19 (identical text/format -> shared string reuse) ---
$ws.Range("E17").Value = "code產出，并通過測試"
$e17chars = $ws.Range("E17").Characters(5, 8)
$e17chars.Font.Name = "宋体"
$e17chars.Font.Size = 12
$e17chars.Font.Color = 16777215

# Copy E17's rich text into E18/E19 so all three rows reuse the same shared
# string entry, same as the source workbook. Use "paste values" so the
# pre-existing cell formatting (fill/border/alignment) of E18/E19 is left
# untouched.
$ws.Range("E17").Copy()
$ws.Range("E18").PasteSpecial(-4163)
$ws.Range("E19").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# --- Column E, rows 15-16 (new accent-font cells) ---------------------------
$ws.Range("E15").Value = "建置資料庫與phpmyadmin中"
$ws.Range("E15").Font.Name = "宋体"
$ws.Range("E15").Font.Size = 12
$ws.Range("E15").Font.Color = 16777215

$ws.Range("E16").Value = "通過測試"
$ws.Range("E16").Font.Name = "宋体"
$ws.Range("E16").Font.Size = 12
$ws.Range("E16").Font.Color = 16777215

# Leave the selection where the author left it when they saved the file.
$ws.Range("E16").Select()
